$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a localization table ("Tabelle2") spanning A1:C46 (Key /
# String EN / String DE). Grow the table by one row via ListRows so the
# table range, autofilter and dimension all expand together, exactly like
# typing a new row into the table in the Excel UI would.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Seed the new row's formatting (fill + wrap text on the two string columns)
# by copying an existing multi-line table row's format, then overwrite the
# values below. Row 14 already carries the "Key plain / String EN+DE
# wrapped" style used throughout the sheet for longer strings.
$ws.Range("A14:C14").Copy()
$ws.Range("A47:C47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A47").Value = "service_worker-update_available"
$ws.Range("B47").Value = "A new update is available. Click the button below to refresh the app and get the latest and greatest stuff!"
$ws.Range("C47").Value = "`nEin neues Update ist verfügbar. Klicken Sie auf die Schaltfläche unten, um die App zu aktualisieren und die neuesten und besten Inhalte zu erhalten!"

# Match the row height used for this wrapped entry.
$ws.Rows.Item(47).RowHeight = 57

# Move the selection to the newly added row, mirroring the author's view.
$ws.Range("A47").Select() | Out-Null
